$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @(2,1,1,0,0,3,1,2,0,2,0,0,1,1,0,0,0,0,0,1,2,1,2,4,2,2,0,1,0,4,2,1,1,0,1,0,1,2,0)

for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
